# Auto-generated-ish PowerShell COM-interop script
# Applies the 'TC for Login added' edit to Test-Cases.xlsx

$wb = $excel.ActiveWorkbook

# ---- 1) Rename sheet 3 and add the new Sheet1 at the end ----
$wsLogin = $wb.Worksheets.Item(3)
$wsLogin.Name = 'Login&Logout_TestCases'

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNotes = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsNotes.Name = 'Sheet1'

# ---- 2) URL_Test_Cases: fill column L (Remarks) with NA ----
$wsUrl = $wb.Worksheets.Item(1)
$wsUrl.Range("L2").Value = 'NA'
$wsUrl.Range("L3").Value = 'NA'
$wsUrl.Range("L4").Value = 'NA'
$wsUrl.Range("L5").Value = 'NA'
$wsUrl.Range("L6").Value = 'NA'
$wsUrl.Range("L2:L6").WrapText = $true
$wsUrl.Range("L2:L6").VerticalAlignment = -4108  # xlCenter

$wsUrl.Range("M4").Select()

# ---- 3) User Registration_TestCases: update J2/K2/F9 ----
$wsReg = $wb.Worksheets.Item(2)
$wsReg.Range("J2").Value = '"Account created successfully '
$wsReg.Range("K2").Value = 'Pass'
$wsReg.Range("F9").Value = 'Phone number should be valid'
$wsReg.Columns.Item(10).ColumnWidth = 17.140625

$wsReg.Range("F10").Select()
$wsReg.Activate()

# ---- 4) Login&Logout_TestCases: populate Test Data / Expected / Actual / Status / Remarks ----
$wsLogin.Columns.Item(9).ColumnWidth = 37.85546875
$wsLogin.Columns.Item(10).ColumnWidth = 19

# Row 3
$wsLogin.Range("H3").Value = 'Valid email, password'
$wsLogin.Range("I3").Value = 'User should be able login to website successfully'
$wsLogin.Range("J3").Value = 'Login Successful'
$wsLogin.Range("K3").Value = 'Pass'
$wsLogin.Range("L3").Value = 'NA'
$wsLogin.Range("H3:L3").WrapText = $true
$wsLogin.Rows.Item(3).RowHeight = 45

# Row 4
$wsLogin.Range("H4").Value = 'Email ,password'
$wsLogin.Range("I4").Value = '"Invalid credentials" or
 "Incorrect email/password"'
$wsLogin.Range("J4").Value = 'Invalid Email/Password
 displayed on UI'
$wsLogin.Range("K4").Value = 'Pass'
$wsLogin.Range("L4").Value = 'NA'
$wsLogin.Range("H4:L4").WrapText = $true
$wsLogin.Rows.Item(4).RowHeight = 60

# Row 5
$wsLogin.Range("H5").Value = 'Email ,password'
$wsLogin.Range("I5").Value = '"Invalid credentials" or
 "Incorrect email/password"'
$wsLogin.Range("J5").Value = 'Invalid Email/Password
 displayed on UI'
$wsLogin.Range("K5").Value = 'Pass'
$wsLogin.Range("L5").Value = 'NA'
$wsLogin.Range("H5:L5").WrapText = $true
$wsLogin.Rows.Item(5).RowHeight = 45

# Row 6
$wsLogin.Range("H6").Value = 'Email ,password'
$wsLogin.Range("I6").Value = '"Invalid credentials" or
 "Incorrect email/password"'
$wsLogin.Range("J6").Value = 'Invalid Email/Password
 displayed on UI'
$wsLogin.Range("K6").Value = 'Pass'
$wsLogin.Range("L6").Value = 'NA'
$wsLogin.Range("H6:L6").WrapText = $true
$wsLogin.Rows.Item(6).RowHeight = 60

# Row 7
$wsLogin.Range("H7").Value = 'Blank data'
$wsLogin.Range("I7").Value = 'Error messages "Email is required" and
 "Password is required" are displayed.'
$wsLogin.Range("J7").Value = 'Email and Password
Required'
$wsLogin.Range("K7").Value = 'Pass'
$wsLogin.Range("L7").Value = 'NA'
$wsLogin.Range("H7:L7").WrapText = $true
$wsLogin.Rows.Item(7).RowHeight = 45

# Row 8
$wsLogin.Range("H8").Value = 'Email ,password'
$wsLogin.Range("I8").Value = '"Invalid credentials" or
 "password incorrect "'
$wsLogin.Range("J8").Value = '"Invalid credentials" or "password incorrect " message displayed when clicked on Login'
$wsLogin.Range("K8").Value = 'Pass'
$wsLogin.Range("L8").Value = 'NA'
$wsLogin.Range("H8:L8").WrapText = $true
$wsLogin.Rows.Item(8).RowHeight = 75

# Row 9
$wsLogin.Range("H9").Value = 'Email ,password'
$wsLogin.Range("I9").Value = '"Password should be masked and clicked on eye should show the password"'
$wsLogin.Range("J9").Value = '"Password is  masked and clicked on eye password is  showen"'
$wsLogin.Range("K9").Value = 'Pass'
$wsLogin.Range("L9").Value = 'NA'
$wsLogin.Range("H9:L9").WrapText = $true
$wsLogin.Rows.Item(9).RowHeight = 60

# Row 10
$wsLogin.Range("H10").Value = 'Email ,password'
$wsLogin.Range("I10").Value = '"Email and Password required "'
$wsLogin.Range("J10").Value = '"Email and Password required "diplayed when clicked on login'
$wsLogin.Range("K10").Value = 'Pass'
$wsLogin.Range("L10").Value = 'NA'
$wsLogin.Range("H10:L10").WrapText = $true
$wsLogin.Rows.Item(10).RowHeight = 60

# Row 11
$wsLogin.Range("H11").Value = 'Email ,password'
$wsLogin.Range("I11").Value = '" Re-directed to Reset password  page"'
$wsLogin.Range("J11").Value = '"Re-directed to password reset page "'
$wsLogin.Range("K11").Value = 'Pass'
$wsLogin.Range("L11").Value = 'NA'
$wsLogin.Range("H11:L11").WrapText = $true
$wsLogin.Rows.Item(11).RowHeight = 60

# Row 12
$wsLogin.Range("H12").Value = 'Proctected URL'
$wsLogin.Range("I12").Value = 'User should be redirected to the login page'
$wsLogin.Range("J12").Value = '"Login page redirected "'
$wsLogin.Range("K12").Value = 'Pass'
$wsLogin.Range("L12").Value = 'NA'
$wsLogin.Range("H12:L12").WrapText = $true
$wsLogin.Rows.Item(12).RowHeight = 135

# Row 13
$wsLogin.Range("I13").Value = 'User should be logged out successfully and
 redirected to the login page or home page'
$wsLogin.Range("J13").Value = '"Login page redirected "'
$wsLogin.Range("K13").Value = 'Pass'
$wsLogin.Range("L13").Value = 'NA'
$wsLogin.Range("H13:L13").WrapText = $true
$wsLogin.Rows.Item(13).RowHeight = 150

# Row 14
$wsLogin.Range("H14").Value = 'Idle time = configured timeout'
$wsLogin.Range("I14").Value = 'User should be logged out automatically 
and redirected to login page |'
$wsLogin.Range("J14").Value = '"user logged out and directed to login page"'
$wsLogin.Range("K14").Value = 'Pass'
$wsLogin.Range("L14").Value = 'NA'
$wsLogin.Range("H14:L14").WrapText = $true
$wsLogin.Rows.Item(14).RowHeight = 90

$wsLogin.Range("E18").Select()

# ---- 5) New Sheet1: add logout-step notes cell ----
$wsNotes.Range("A1").Value = '
1.Open the application
2.Ensure the user is logged in
3.Click on Profile / Account / Logout option
4.Confirm logout (if confirmation popup is shown)
5.Observe the system behavior
'
$wsNotes.Range("A1").WrapText = $true
$wsNotes.Columns.Item(1).ColumnWidth = 49.85546875
$wsNotes.Rows.Item(1).RowHeight = 255

# ---- 6) Make User Registration_TestCases the active tab, matching the saved view ----
$wsReg.Activate()
$wsReg.Range("F10").Select()
